# Apply the "Managed to make a variable endpoint" edit to the
# GhostsGoblinsEnemiesPlan workbook.
#
# Summary of the change:
#  * A new column is inserted before column G, shifting the second
#    "day-of-week" header block (old G:K) one column to the right (H:L),
#    and duplicating the "Friday" header into the new column G.
#  * A center+wrap-without-fill style is used for a couple of new note
#    cells, and an existing note cell (E3) picks up that same alignment.
#  * New task/notes text is added in E4, E5, D6 and E6.
#  * Row 6's height grows to 60 to fit the new wrapped text.
#  * The active selection moves to D7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before G so the second day-header block (which used
# to start at G) slides right to H:L, leaving room for a duplicated
# "Friday" header in the new column G.
$ws.Range("G1").EntireColumn.Insert()

# Populate the newly freed column G with "Friday", matching B1.
$ws.Range("G1").Value2 = $ws.Range("B1").Value2

# Add the new note / task text. E4 just gets typed into (keeping its
# existing, unformatted style) - no alignment change there.
$ws.Range("E4").Value2 = "Pseudocode Zombie Spawn"
$ws.Range("E5").Value2 = "Pseudocode Zombie Spawner Behavior"
$ws.Range("E6").Value2 = "Implement Sprite and figure out how to reverse zombie movement in spawner"
$ws.Range("D6").Value2 = "Make target a variable off camera"

# Center + wrap text was applied across the new note cells (E3, E5, D6,
# E6).
foreach ($addr in @("E3", "E5", "D6", "E6")) {
    $cell = $ws.Range($addr)
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.WrapText = $true
}

# Row 6 grows taller to accommodate the new wrapped text.
$ws.Range("A6").RowHeight = 60

# Move the active selection, as captured in the saved workbook.
$ws.Range("D7").Select()
